$wb = $excel.ActiveWorkbook

# --- Data sheet: clear the stray duplicate "Aust" label in N3 ---
$wsData = $wb.Worksheets.Item("Data")
$wsData.Range("N3").Value = $null

# --- Progress sheet: new widget def needs two blank leading columns
#     (B:C) ahead of the state columns, so the whole state block shifts
#     from B:J to D:L. Remove the now-redundant row-2 label too. ---
$wsProgress = $wb.Worksheets.Item("Progress")
$wsProgress.Range("B1:C1").EntireColumn.Insert()
$wsProgress.Range("A2").Value = $null

# --- Description sheet: status wording update ---
$wsDesc = $wb.Worksheets.Item("Description")
$wsDesc.Range("B2").Value = "On track"

# --- View state: Description tab is now the active/selected tab, with
#     updated selections left on the other two sheets ---
$wsData.Range("M17").Select()
$wsProgress.Range("B1").Select()
$wsDesc.Activate()
$wsDesc.Range("B3").Select()

Write-Output "done"
